# Applies the "parte 1 de nuevos estado de cuenta" update:
#  - swaps the "Periodo Mora" values shown for the two rows of
#    NIBIA ESTHER MEDRANO PATERNINA (1805 <-> 1804)
#  - updates the corresponding "Valor Mora" amounts from 781242 to 737717

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Periodo Mora (stored as text, numFmt "@") - swap the two period labels
$ws.Range("E16").Value = "1804"
$ws.Range("E17").Value = "1805"

# Valor Mora - updated amounts
$ws.Range("G16").Value = 737717
$ws.Range("G17").Value = 737717
